$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link  = "https://www.360dx.com/cancer/fda-approves-agilent-cdx-assay-use-bristol-myers-squibb-immunotherapies"
$kw    = "CDx"
$title = "FDA Approves Agilent CDx Assay for Use With Bristol Myers Squibb Immunotherapies"

$row = 36

# Write the new row's values first (link text doubles as the hyperlink display text,
# matching the pattern used by every other row in the sheet).
$ws.Cells.Item($row, 1).Value = $link
$ws.Cells.Item($row, 2).Value = $kw
$ws.Cells.Item($row, 3).Value = $title

# Attach the hyperlink relationship for the new cell.
$ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $link)

# Hyperlinks.Add stamps a freshly minted cell style; put the cell back on the
# same shared "Hyperlink" style every other link cell in column A already uses.
$ws.Cells.Item($row, 1).Style = "Hyperlink"
